# Updated cryptos list on Sun Feb  4 22:32:26 UTC 2024 with GitHub Actions
# Applies Price (column D) and Volume(1h) (column E) updates for rows 2-51.
#
# Note: values are prefixed with a leading apostrophe so Excel stores them as
# literal text (matching the workbook's original inline-string cells, e.g.
# "306.10", "32.80", "42.640.50") instead of auto-converting them to numbers,
# which would drop trailing zeros / change their formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.640.50"
$ws.Range("E2").Value = "'  -0.99%  "
$ws.Range("D3").Value = "'2.286.48"
$ws.Range("E3").Value = "'  -0.55%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'306.10"
$ws.Range("E5").Value = "'  +1.81%  "
$ws.Range("D6").Value = "'95.73"
$ws.Range("E6").Value = "'  -2.63%  "
$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "'  -2.62%  "
$ws.Range("E8").Value = "'  +0.11%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "'  -3.05%  "
$ws.Range("D10").Value = "'35.04"
$ws.Range("E10").Value = "'  -3.19%  "
$ws.Range("D11").Value = "'0.0787"
$ws.Range("E11").Value = "'  -0.36%  "
$ws.Range("D12").Value = "'18.37"
$ws.Range("E12").Value = "'  +3.68%  "
$ws.Range("E13").Value = "'  +1.12%  "
$ws.Range("D14").Value = "'6.72"
$ws.Range("E14").Value = "'  -2.15%  "
$ws.Range("D15").Value = "'2.654.93"
$ws.Range("E15").Value = "'  -0.07%  "
$ws.Range("D16").Value = "'2.290.50"
$ws.Range("E16").Value = "'  -0.20%  "
$ws.Range("D17").Value = "'0.778"
$ws.Range("E17").Value = "'  -1.30%  "
$ws.Range("D18").Value = "'42.615.56"
$ws.Range("E18").Value = "'  -0.76%  "
$ws.Range("D19").Value = "'12.86"
$ws.Range("E19").Value = "'  +0.24%  "
$ws.Range("D20").Value = "'0.0₃0893"
$ws.Range("E20").Value = "'  -2.14%  "
$ws.Range("D21").Value = "'6.04"
$ws.Range("E21").Value = "'  -1.42%  "
$ws.Range("D22").Value = "'66.99"
$ws.Range("E22").Value = "'  -2.85%  "
$ws.Range("D23").Value = "'235.03"
$ws.Range("E23").Value = "'  -0.97%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "'  -0.41%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "'  +0.76%  "
$ws.Range("E26").Value = "'  +0.07%  "
$ws.Range("E27").Value = "'  +0.03%  "
$ws.Range("D28").Value = "'24.98"
$ws.Range("E28").Value = "'  +0.14%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = "'  +6.79%  "
$ws.Range("D30").Value = "'166.11"
$ws.Range("E30").Value = "'  +0.85%  "
$ws.Range("D31").Value = "'9.03"
$ws.Range("E31").Value = "'  -0.96%  "
$ws.Range("D32").Value = "'32.80"
$ws.Range("E32").Value = "'  -0.75%  "
$ws.Range("E33").Value = "'  +0.10%  "
$ws.Range("D34").Value = "'4.73"
$ws.Range("E34").Value = "'  -0.81%  "
$ws.Range("D35").Value = "'4.96"
$ws.Range("E35").Value = "'  -2.35%  "
$ws.Range("D36").Value = "'17.52"
$ws.Range("E36").Value = "'  -2.26%  "
$ws.Range("E37").Value = "'  -0.39%  "
$ws.Range("D38").Value = "'0.0691"
$ws.Range("E38").Value = "'  -0.79%  "
$ws.Range("E39").Value = "'  -1.06%  "
$ws.Range("D40").Value = "'1.74"
$ws.Range("E40").Value = "'  -2.12%  "
$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "'  -1.71%  "
$ws.Range("D42").Value = "'2.69"
$ws.Range("E42").Value = "'  -3.57%  "
$ws.Range("D43").Value = "'1.999.70"
$ws.Range("E43").Value = "'  -0.64%  "
$ws.Range("D44").Value = "'0.0278"
$ws.Range("E44").Value = "'  -2.83%  "
$ws.Range("D45").Value = "'18.13"
$ws.Range("E45").Value = "'  +3.78%  "
$ws.Range("D46").Value = "'10.04"
$ws.Range("E46").Value = "'  -2.82%  "
$ws.Range("D47").Value = "'2.02"
$ws.Range("E47").Value = "'  -9.46%  "
$ws.Range("D48").Value = "'2.78"
$ws.Range("E48").Value = "'  -1.57%  "
$ws.Range("D49").Value = "'2.91"
$ws.Range("E49").Value = "'  +9.54%  "
$ws.Range("D50").Value = "'53.72"
$ws.Range("E50").Value = "'  -0.65%  "
$ws.Range("D51").Value = "'2.519.24"
$ws.Range("E51").Value = "'  -0.17%  "
